# Update ARC_YR_FIN.xlsx: add the latest fiscal-period column (new column D)
# in front of the existing year-over-year data, shifting the historical
# columns (old D:K) one place to the right (new E:L), then populate the
# new column with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D; this shifts D:K -> E:L
# (including formulas-free literal values) and grows the sheet dimension.
$ws.Columns("D:D").Insert()

# The freshly inserted column D has no number format yet; copy the
# formatting (date format row 7/38/80, numeric format elsewhere) from the
# column immediately to its right (old D, now E) so every cell in D keeps
# the same look as the rest of its row. Only touch the row ranges that
# actually carry D:K data (skip the blank separator rows 36, 37, 78, 79
# which have no cells there, so we don't create stray empty D cells).
$ws.Range("E7:E35").Copy() | Out-Null
$ws.Range("D7:D35").PasteSpecial(-4122) | Out-Null

$ws.Range("E38:E77").Copy() | Out-Null
$ws.Range("D38:D77").PasteSpecial(-4122) | Out-Null

$ws.Range("E80:E102").Copy() | Out-Null
$ws.Range("D80:D102").PasteSpecial(-4122) | Out-Null

# ---- Income Statement (rows 7-35) ----
$ws.Range("D7").Value = 43465

$ws.Range("D8").Value = 400800
$ws.Range("D9").Value = 269900
$ws.Range("D10").Value = 130900
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 3900

$ws.Range("D17").Value = 382900
$ws.Range("D18").Value = 17900

$ws.Range("D20").Value = -5800
$ws.Range("D21").Value = 44900
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 12100
$ws.Range("D24").Value = 3300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 8700
$ws.Range("D27").Value = 8900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 5800
$ws.Range("D33").Value = 8900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 8900

# ---- Balance Sheet (rows 38-77) ----
$ws.Range("D38").Value = 43465

$ws.Range("D41").Value = 29400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 58000
$ws.Range("D44").Value = 16800
$ws.Range("D45").Value = 11100
$ws.Range("D46").Value = 115400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 70700
$ws.Range("D49").Value = 126200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 27500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 339700

$ws.Range("D57").Value = 24200
$ws.Range("D58").Value = 22100
$ws.Range("D59").Value = 34600
$ws.Range("D60").Value = 81000
$ws.Range("D61").Value = 105100
$ws.Range("D62").Value = 6400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 199400

$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 29400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 140300
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (rows 80-102) ----
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 8900

$ws.Range("D83").Value = 32900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 55000

$ws.Range("D91").Value = -14900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -14200

$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -38700
$ws.Range("D101").Value = -700
$ws.Range("D102").Value = 1400
